# Penalty Reward System (unfinished) - adjust weekly/monthly PO data
$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Weekly Quantity" ---
$ws1 = $wb.Worksheets.Item("Weekly Quantity")

# Remove row 5 entirely (shifts rows 6-26 up to 5-25)
$ws1.Rows.Item(5).Delete()

# Update requested quantity for row 4 (week of 45361.99999999999)
$ws1.Range("B4").Value = 60

# --- Sheet 2: "Monthly Trend" ---
$ws2 = $wb.Worksheets.Item("Monthly Trend")

# Update requested quantity for row 3 (month of 45382.99999999999)
$ws2.Range("B3").Value = 60
